$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the stray leading spaces from the header/text cells
$ws.Range("B1").Value = "trow"
$ws.Range("C1").Value = "frow"
$ws.Range("B2").Value = "text 1"
$ws.Range("B3").Value = "text 2"
$ws.Range("B4").Value = "text 3"

# Update the saved selection in the sheet view
$ws.Range("B5").Select()
